$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# ERT_ATFM_YY sheet: release date + row 15 raw input
# -----------------------------------------------------------------
$wsYY = $wb.Worksheets.Item("ERT_ATFM_YY")
$wsYY.Range("B2").Value = 45758
$wsYY.Range("E15").Value = 20875950

# -----------------------------------------------------------------
# ERT_ATFM_MM sheet: monthly "FLTS" delay minutes (col D), rows 66-77
# -----------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("ERT_ATFM_MM")
$wsMM.Range("D66").Value = 283577
$wsMM.Range("D67").Value = 228949
$wsMM.Range("D68").Value = 387745
$wsMM.Range("D69").Value = 465862
$wsMM.Range("D70").Value = 1455200
$wsMM.Range("D71").Value = 3667211
$wsMM.Range("D72").Value = 5577411
$wsMM.Range("D73").Value = 4141118
$wsMM.Range("D74").Value = 2320537
$wsMM.Range("D75").Value = 1325196
$wsMM.Range("D76").Value = 485228
$wsMM.Range("D77").Value = 537916

# -----------------------------------------------------------------
# ERT_ATFM_FAB sheet: FAB-level delay minutes (col D), rows 6-15
# -----------------------------------------------------------------
$wsFAB = $wb.Worksheets.Item("ERT_ATFM_FAB")
$wsFAB.Range("D6").Value = 20875950
$wsFAB.Range("D7").Value = 171611
$wsFAB.Range("D8").Value = 2503737
$wsFAB.Range("D9").Value = 206910
$wsFAB.Range("D11").Value = 5176583
$wsFAB.Range("D12").Value = 10068554
$wsFAB.Range("D13").Value = 26033
$wsFAB.Range("D14").Value = 2686304
$wsFAB.Range("D15").Value = 333036

# -----------------------------------------------------------------
# ERT_ATFM_LOC sheet: Entity table, rows 6-33 renamed / renumbered,
# columns B/E/F converted from formulas to plain values, row 34
# cleared (table shrinks from 29 to 28 entities).
# -----------------------------------------------------------------
$wsLOC = $wb.Worksheets.Item("ERT_ATFM_LOC")

$wsLOC.Range("A6").Value = "ANS CR"
$wsLOC.Range("B6").Value = 0.11
$wsLOC.Range("C6").Value = 711214
$wsLOC.Range("D6").Value = 89658
$wsLOC.Range("E6").Value = 0.13
$wsLOC.Range("F6").Value = 0.02

$wsLOC.Range("A7").Value = "Austro Control"
$wsLOC.Range("B7").Value = 0.16
$wsLOC.Range("C7").Value = 1249634
$wsLOC.Range("D7").Value = 594968
$wsLOC.Range("E7").Value = 0.48
$wsLOC.Range("F7").Value = 0.32

$wsLOC.Range("A8").Value = "Avinor Flysikring AS"
$wsLOC.Range("B8").Value = 0.11
$wsLOC.Range("C8").Value = 555589
$wsLOC.Range("D8").Value = 24554
$wsLOC.Range("E8").Value = 0.04
$wsLOC.Range("F8").Value = -0.07

$wsLOC.Range("A9").Value = "BULATSA"
$wsLOC.Range("B9").Value = 0.08
$wsLOC.Range("C9").Value = 1050462
$wsLOC.Range("D9").Value = 102283
$wsLOC.Range("E9").Value = 0.1
$wsLOC.Range("F9").Value = 0.02

$wsLOC.Range("A10").Value = "Croatia Control"
$wsLOC.Range("B10").Value = 0.17
$wsLOC.Range("C10").Value = 920123
$wsLOC.Range("D10").Value = 1389663
$wsLOC.Range("E10").Value = 1.51
$wsLOC.Range("F10").Value = 1.34

$wsLOC.Range("A11").Value = "DCAC Cyprus"
$wsLOC.Range("B11").Value = 0.15
$wsLOC.Range("C11").Value = 379007
$wsLOC.Range("D11").Value = 2061
$wsLOC.Range("E11").Value = 0.01
$wsLOC.Range("F11").Value = -0.14

$wsLOC.Range("A12").Value = "DFS + MUAC-DE"
$wsLOC.Range("B12").Value = 0.27
$wsLOC.Range("C12").Value = 2839817
$wsLOC.Range("D12").Value = 4467388
$wsLOC.Range("E12").Value = 1.57
$wsLOC.Range("F12").Value = 1.3

$wsLOC.Range("A13").Value = "DSNA"
$wsLOC.Range("B13").Value = 0.25
$wsLOC.Range("C13").Value = 3391028
$wsLOC.Range("D13").Value = 4708979
$wsLOC.Range("E13").Value = 1.39
$wsLOC.Range("F13").Value = 1.14

$wsLOC.Range("A14").Value = "EANS"
$wsLOC.Range("B14").Value = 0.03
$wsLOC.Range("C14").Value = 169709
$wsLOC.Range("D14").Value = 1132
$wsLOC.Range("E14").Value = 0.01
$wsLOC.Range("F14").Value = -0.02

$wsLOC.Range("A15").Value = "ENAIRE"
$wsLOC.Range("B15").Value = 0.19
$wsLOC.Range("C15").Value = 2360451
$wsLOC.Range("D15").Value = 2403805
$wsLOC.Range("E15").Value = 1.02
$wsLOC.Range("F15").Value = 0.83

$wsLOC.Range("A16").Value = "ENAV"
$wsLOC.Range("B16").Value = 0.11
$wsLOC.Range("C16").Value = 2027066
$wsLOC.Range("D16").Value = 1454304
$wsLOC.Range("E16").Value = 0.72
$wsLOC.Range("F16").Value = 0.61

$wsLOC.Range("A17").Value = "Fintraffic ANS"
$wsLOC.Range("B17").Value = 0.05
$wsLOC.Range("C17").Value = 230292
$wsLOC.Range("D17").Value = 0
$wsLOC.Range("E17").Value = 0
$wsLOC.Range("F17").Value = -0.05

$wsLOC.Range("A18").Value = "HASP"
$wsLOC.Range("B18").Value = 0.15
$wsLOC.Range("C18").Value = 1074741
$wsLOC.Range("D18").Value = 1047372
$wsLOC.Range("E18").Value = 0.97
$wsLOC.Range("F18").Value = 0.82

$wsLOC.Range("A19").Value = "HungaroControl (EC)"
$wsLOC.Range("B19").Value = 0.11
$wsLOC.Range("C19").Value = 1095404
$wsLOC.Range("D19").Value = 3014149
$wsLOC.Range("E19").Value = 2.75
$wsLOC.Range("F19").Value = 2.64

$wsLOC.Range("A20").Value = "AirNav Ireland"
$wsLOC.Range("B20").Value = 0.03
$wsLOC.Range("C20").Value = 679303
$wsLOC.Range("D20").Value = 1531
$wsLOC.Range("E20").Value = 0
$wsLOC.Range("F20").Value = -0.03

$wsLOC.Range("A21").Value = "LFV"
$wsLOC.Range("B21").Value = 0.08
$wsLOC.Range("C21").Value = 605374
$wsLOC.Range("D21").Value = 6357
$wsLOC.Range("E21").Value = 0.01
$wsLOC.Range("F21").Value = -0.07

$wsLOC.Range("A22").Value = "LGS"
$wsLOC.Range("B22").Value = 0.03
$wsLOC.Range("C22").Value = 227609
$wsLOC.Range("D22").Value = 347
$wsLOC.Range("E22").Value = 0
$wsLOC.Range("F22").Value = -0.03

$wsLOC.Range("A23").Value = "LPS SR"
$wsLOC.Range("B23").Value = 0.07
$wsLOC.Range("C23").Value = 601041
$wsLOC.Range("D23").Value = 66932
$wsLOC.Range("E23").Value = 0.11
$wsLOC.Range("F23").Value = 0.04

$wsLOC.Range("A24").Value = "LVNL + MUAC-NL"
$wsLOC.Range("B24").Value = 0.14
$wsLOC.Range("C24").Value = 1201768
$wsLOC.Range("D24").Value = 114786
$wsLOC.Range("E24").Value = 0.1
$wsLOC.Range("F24").Value = -0.04

$wsLOC.Range("A25").Value = "Malta Air Traffic Services Ltd."
$wsLOC.Range("B25").Value = 0.01
$wsLOC.Range("C25").Value = 151700
$wsLOC.Range("D25").Value = 0
$wsLOC.Range("E25").Value = 0
$wsLOC.Range("F25").Value = -0.01

$wsLOC.Range("A26").Value = "NAV Portugal (Continental)"
$wsLOC.Range("B26").Value = 0.13
$wsLOC.Range("C26").Value = 720113
$wsLOC.Range("D26").Value = 282499
$wsLOC.Range("E26").Value = 0.39
$wsLOC.Range("F26").Value = 0.26

$wsLOC.Range("A27").Value = "NAVIAIR"
$wsLOC.Range("B27").Value = 0.05
$wsLOC.Range("C27").Value = 605140
$wsLOC.Range("D27").Value = 28330
$wsLOC.Range("E27").Value = 0.05
$wsLOC.Range("F27").Value = 0

$wsLOC.Range("A28").Value = "PANSA"
$wsLOC.Range("B28").Value = 0.12
$wsLOC.Range("C28").Value = 745435
$wsLOC.Range("D28").Value = 171611
$wsLOC.Range("E28").Value = 0.23
$wsLOC.Range("F28").Value = 0.11

$wsLOC.Range("A29").Value = "ROMATSA"
$wsLOC.Range("B29").Value = 0.04
$wsLOC.Range("C29").Value = 819673
$wsLOC.Range("D29").Value = 104627
$wsLOC.Range("E29").Value = 0.13
$wsLOC.Range("F29").Value = 0.09

$wsLOC.Range("A30").Value = "SE Oro Navigacija"
$wsLOC.Range("B30").Value = 0.02
$wsLOC.Range("C30").Value = 187523
$wsLOC.Range("D30").Value = 0
$wsLOC.Range("E30").Value = 0
$wsLOC.Range("F30").Value = -0.02

$wsLOC.Range("A31").Value = "Belgium-Lux. + MUAC  BE-LU"
$wsLOC.Range("B31").Value = 0.17
$wsLOC.Range("C31").Value = 1215122
$wsLOC.Range("D31").Value = 237678
$wsLOC.Range("E31").Value = 0.2
$wsLOC.Range("F31").Value = 0.03

$wsLOC.Range("A32").Value = "Skyguide"
$wsLOC.Range("B32").Value = 0.19
$wsLOC.Range("C32").Value = 1334046
$wsLOC.Range("D32").Value = 539723
$wsLOC.Range("E32").Value = 0.4
$wsLOC.Range("F32").Value = 0.21

$wsLOC.Range("A33").Value = "Slovenia Control, Ltd"
$wsLOC.Range("B33").Value = 0.09
$wsLOC.Range("C33").Value = 429597
$wsLOC.Range("D33").Value = 21213
$wsLOC.Range("E33").Value = 0.05
$wsLOC.Range("F33").Value = -0.04

# Row 34 previously held "Slovenia Control" data; the table now only
# has 28 entities (rows 6-33), so row 34 is cleared out completely and
# restyled to match the blank trailer rows below it (e.g. row 35).
$rowSrc = $wsLOC.Range("A35:F35")
$rowDst = $wsLOC.Range("A34:F34")
$rowSrc.Copy()
$rowDst.PasteSpecial(-4122)
$rowDst.ClearContents()

Write-Host "Edit script completed"
